$wb = $excel.ActiveWorkbook

# --- Update the ProjectModule sheet view (selection moves from J16 to G25) ---
$wsProject = $wb.Worksheets.Item("ProjectModule")
$wsProject.Activate()
$wsProject.Range("G25").Select()

# --- Update the Submission sheet view (selection moves from I12 to L15) ---
$wsSubmission = $wb.Worksheets.Item("Submission")
$wsSubmission.Activate()

# Update submissionDueDate column (E) text values to include the time component
for ($r = 2; $r -le 19; $r++) {
    $wsSubmission.Range("E$r").Value = "2024-08-01 00:00:00"
}
for ($r = 20; $r -le 29; $r++) {
    $wsSubmission.Range("E$r").Value = "2022-04-09 00:00:00"
}

$wsSubmission.Range("L15").Select()
